$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 126.333336
$ws.Range("I4").Value = 126.333336
$ws.Range("K4").Value = 126.333336
$ws.Range("M4").Value = -12.333336

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1342.5714
$ws.Range("J29").Value = 899.6667
$ws.Range("L29").Value = 2699.0001
$ws.Range("N29").Value = -3261.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4939.3
$ws.Range("I76").Value = 4998.6665
$ws.Range("K76").Value = 4998.6665
$ws.Range("M76").Value = -4683.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4939.3
$ws.Range("I79").Value = 4998.6665
$ws.Range("K79").Value = 4998.6665
$ws.Range("M79").Value = -3906.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1542009.5
$ws.Range("I86").Value = 2940870.8
$ws.Range("J86").Value = 3262.1
$ws.Range("K86").Value = 2940870.8
$ws.Range("L86").Value = 3262.1
$ws.Range("M86").Value = -2939747.8
$ws.Range("N86").Value = -5508.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1542009.5
$ws.Range("I89").Value = 2940870.8
$ws.Range("J89").Value = 3262.1
$ws.Range("K89").Value = 14704354
$ws.Range("L89").Value = 16310.5
$ws.Range("M89").Value = -14698738
$ws.Range("N89").Value = -27542.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 672.75
$ws.Range("I92").Value = 590.6
$ws.Range("K92").Value = 590.6
$ws.Range("M92").Value = 657.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2757
$ws.Range("J125").Value = 3883
$ws.Range("L125").Value = 34947
$ws.Range("N125").Value = -39867

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 12177.846
$ws.Range("I132").Value = 5283.9375
$ws.Range("J132").Value = 15241.806
$ws.Range("K132").Value = 15851.8125
$ws.Range("L132").Value = 45725.41800000001
$ws.Range("M132").Value = -13321.8125
$ws.Range("N132").Value = -50785.41800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 608.1667
$ws.Range("I5").Value = 608.1667
$ws.Range("K5").Value = 608.1667
$ws.Range("M5").Value = -496.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8066081
$ws.Range("I74").Value = 13889957
$ws.Range("K74").Value = 13889957
$ws.Range("M74").Value = -13889083

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 8066081
$ws.Range("I77").Value = 13889957
$ws.Range("K77").Value = 69449785
$ws.Range("M77").Value = -69445417

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 109931.336
$ws.Range("J109").Value = 109931.336
$ws.Range("L109").Value = 109931.336
$ws.Range("N109").Value = -112705.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 82499.5
$ws.Range("J112").Value = 82499.5
$ws.Range("L112").Value = 82499.5
$ws.Range("N112").Value = -85453.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2153.4
$ws.Range("I122").Value = 2222.8235
$ws.Range("K122").Value = 6668.470499999999
$ws.Range("M122").Value = -4218.470499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 60342.273
$ws.Range("I132").Value = 95761
$ws.Range("K132").Value = 287283
$ws.Range("M132").Value = -284753

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 608.1667
$ws.Range("I4").Value = 608.1667
$ws.Range("K4").Value = 608.1667
$ws.Range("M4").Value = -493.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2652
$ws.Range("I20").Value = 3053.6365
$ws.Range("J20").Value = 1179.3334
$ws.Range("K20").Value = 3053.6365
$ws.Range("L20").Value = 1179.3334
$ws.Range("M20").Value = -2806.6365
$ws.Range("N20").Value = -1673.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 586.7143
$ws.Range("I22").Value = 586.7143
$ws.Range("K22").Value = 586.7143
$ws.Range("M22").Value = -413.7143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2923.1538
$ws.Range("I105").Value = 2666.0667
$ws.Range("K105").Value = 2666.0667
$ws.Range("M105").Value = -919.0666999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 87750
$ws.Range("J112").Value = 87750
$ws.Range("L112").Value = 87750
$ws.Range("N112").Value = -90704

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 246.13333
$ws.Range("I7").Value = 200
$ws.Range("J7").Value = 257.66666
$ws.Range("K7").Value = 200
$ws.Range("L7").Value = 257.66666
$ws.Range("M7").Value = -87
$ws.Range("N7").Value = -483.66666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1799.0454
$ws.Range("I16").Value = 1709.5
$ws.Range("K16").Value = 1709.5
$ws.Range("M16").Value = -1422.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 620.7273
$ws.Range("I22").Value = 311.14285
$ws.Range("J22").Value = 1162.5
$ws.Range("K22").Value = 311.14285
$ws.Range("L22").Value = 1162.5
$ws.Range("M22").Value = 38.85714999999999
$ws.Range("N22").Value = -1862.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1264124.1
$ws.Range("I105").Value = 1895253.4
$ws.Range("K105").Value = 1895253.4
$ws.Range("M105").Value = -1893506.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1799.0454
$ws.Range("I113").Value = 1709.5
$ws.Range("K113").Value = 1709.5
$ws.Range("M113").Value = 460.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1874.5385
$ws.Range("I122").Value = 1758.3529
$ws.Range("K122").Value = 5275.0587
$ws.Range("M122").Value = -2825.0587

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 10103190
$ws.Range("I132").Value = 13335366
$ws.Range("J132").Value = 2639.125
$ws.Range("K132").Value = 40006098
$ws.Range("L132").Value = 7917.375
$ws.Range("M132").Value = -40003568
$ws.Range("N132").Value = -12977.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 478.16666
$ws.Range("I12").Value = 415
$ws.Range("J12").Value = 490.8
$ws.Range("K12").Value = 1245
$ws.Range("L12").Value = 1472.4
$ws.Range("M12").Value = -1072
$ws.Range("N12").Value = -1818.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5510.846
$ws.Range("I80").Value = 5448
$ws.Range("J80").Value = 5522.273
$ws.Range("K80").Value = 16344
$ws.Range("L80").Value = 16566.819
$ws.Range("M80").Value = -15408
$ws.Range("N80").Value = -18438.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 5510.846
$ws.Range("I83").Value = 5448
$ws.Range("J83").Value = 5522.273
$ws.Range("K83").Value = 49032
$ws.Range("L83").Value = 49700.457
$ws.Range("M83").Value = -44352
$ws.Range("N83").Value = -59060.457

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 684
$ws.Range("I122").Value = 699
$ws.Range("J122").Value = 678
$ws.Range("K122").Value = 6291
$ws.Range("L122").Value = 6102
$ws.Range("M122").Value = -3841
$ws.Range("N122").Value = -11002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1671.8889
$ws.Range("J129").Value = 1529.1666
$ws.Range("L129").Value = 4587.4998
$ws.Range("N129").Value = -14587.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 7872545
$ws.Range("J137").Value = 8477741
$ws.Range("L137").Value = 25433223
$ws.Range("N137").Value = -25443423

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 720476.8
$ws.Range("I80").Value = 1425053.6
$ws.Range("J80").Value = 15900
$ws.Range("K80").Value = 1425053.6
$ws.Range("L80").Value = 15900
$ws.Range("M80").Value = -1424055.6
$ws.Range("N80").Value = -17896

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 720476.8
$ws.Range("I83").Value = 1425053.6
$ws.Range("J83").Value = 15900
$ws.Range("K83").Value = 7125268
$ws.Range("L83").Value = 79500
$ws.Range("M83").Value = -7120276
$ws.Range("N83").Value = -89484

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1999.6
$ws.Range("I122").Value = 1999.6
$ws.Range("K122").Value = 5998.799999999999
$ws.Range("M122").Value = -3548.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 905788.7
$ws.Range("J134").Value = 905788.7
$ws.Range("L134").Value = 2717366.1
$ws.Range("N134").Value = -2722436.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4665.778
$ws.Range("I22").Value = 2001
$ws.Range("K22").Value = 2001
$ws.Range("M22").Value = -1706

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 4665.778
$ws.Range("I27").Value = 2001
$ws.Range("K27").Value = 2001
$ws.Range("M27").Value = -1894

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10486658
$ws.Range("J132").Value = 16673057
$ws.Range("L132").Value = 50019171
$ws.Range("N132").Value = -50024231
